$wb = $excel.ActiveWorkbook

# --- 1) Update the "总计" (summary) sheet: insert the 2022-Q3 row, shift the rest down ---
$wsTotal = $wb.Worksheets.Item(1)

# Row 2 becomes the new 2022-Q3 entry
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.09

# Row 3 becomes what used to be row 2 (2022-Q1)
$wsTotal.Range("B3").Value = "2022-Q1"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.05

# Row 4 is new and becomes what used to be row 3 (2021-Q3); copy A3's style for the index cell
$wsTotal.Range("A3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q3"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.06

# --- 2) Insert a brand-new worksheet "2022-Q3" right after "总计", before the existing "2022-Q1" sheet ---
$newWs = $wb.Worksheets.Add($null, $wsTotal)
$newWs.Name = "2022-Q3"

# Re-fetch the existing "2022-Q1" sheet (it has shifted to position 3 after the insert)
$wsQ1 = $wb.Worksheets.Item(3)

# Clone header row + row formatting (style) from the "2022-Q1" fund-holdings sheet
$wsQ1.Range("B1:H1").Copy($newWs.Range("B1:H1"))
$wsQ1.Range("A2").Copy($newWs.Range("A2"))
$wsQ1.Range("A2").Copy($newWs.Range("A3"))

# Row 2: fund 014232 (text columns use a leading apostrophe so codes/figures are stored as text,
# matching how fund code / figure columns are stored elsewhere in the workbook)
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "'014232"
$newWs.Range("C2").Value = "博时专精特新主题混合A"
$newWs.Range("D2").Value = "'3.14"
$newWs.Range("E2").Value = "'81.93"
$newWs.Range("F2").Value = "'1.62"
$newWs.Range("G2").Value = "'0.0509"
$newWs.Range("H2").Value = 6

# Row 3: fund 014233
$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "'014233"
$newWs.Range("C3").Value = "博时专精特新主题混合C"
$newWs.Range("D3").Value = "'2.69"
$newWs.Range("E3").Value = "'81.93"
$newWs.Range("F3").Value = "'1.62"
$newWs.Range("G3").Value = "'0.0436"
$newWs.Range("H3").Value = 6
